$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font, border, centered alignment) used by the
# existing label column cells (A5:A6, style index 1) down onto the two new
# label cells before filling in their values.
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)

# Row 7 ("a1")
$ws.Range("A7").Value2 = "a1"
$ws.Range("B7").Value2 = 0.8646729588508606
$ws.Range("C7").Value2 = 0.4265280067920685
$ws.Range("D7").Value2 = 0.8276968598365784
$ws.Range("E7").Value2 = 0.8786906003952026
$ws.Range("F7").Value2 = 0.8721588850021362
$ws.Range("G7").Value2 = 97.02021026611328
$ws.Range("H7").Value2 = 12.51860427856445
$ws.Range("I7").Value2 = 12.0368595123291
$ws.Range("J7").Value2 = 178.1204681396484
$ws.Range("K7").Value2 = 194.4346466064453

# Row 8 ("b2")
$ws.Range("A8").Value2 = "b2"
$ws.Range("B8").Value2 = 0.8646729588508606
$ws.Range("C8").Value2 = 0.4265280067920685
$ws.Range("D8").Value2 = 0.8276968598365784
$ws.Range("E8").Value2 = 0.8786906003952026
$ws.Range("F8").Value2 = 0.8721616268157959
$ws.Range("G8").Value2 = 97.02021026611328
$ws.Range("H8").Value2 = 12.52402591705322
$ws.Range("I8").Value2 = 12.03680610656738
$ws.Range("J8").Value2 = 178.1204681396484
$ws.Range("K8").Value2 = 194.4124450683594
